# Updated cryptos list values (prices and volume %) per diff.
# Using an apostrophe prefix on assignment forces Excel to keep these as
# literal text (matching the source inlineStr cells) instead of
# auto-converting numeric-looking strings (e.g. "1.00", "0.981") into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'96.358.95"
$ws.Range("E2").Value = "'  -1.17%  "
$ws.Range("D3").Value = "'3.323.02"
$ws.Range("E3").Value = "'  -2.58%  "
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("D5").Value = "'248.28"
$ws.Range("E5").Value = "'  -2.50%  "
$ws.Range("D6").Value = "'650.54"
$ws.Range("E6").Value = "'  -0.29%  "
$ws.Range("E7").Value = "'  -6.45%  "
$ws.Range("E8").Value = "'  -1.11%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "'  +0.03%  "
$ws.Range("D10").Value = "'0.981"
$ws.Range("E10").Value = "'  -7.12%  "
$ws.Range("D11").Value = "'3.320.32"
$ws.Range("E11").Value = "'  -2.55%  "
$ws.Range("E12").Value = "'  -3.41%  "
$ws.Range("D13").Value = "'40.02"
$ws.Range("E13").Value = "'  -4.18%  "
$ws.Range("D14").Value = "'96.066.94"
$ws.Range("E14").Value = "'  -1.15%  "
$ws.Range("D15").Value = "'6.04"
$ws.Range("E15").Value = "'  -3.63%  "
$ws.Range("E16").Value = "'  -3.89%  "
$ws.Range("D17").Value = "'3.936.95"
$ws.Range("E17").Value = "'  -2.57%  "
$ws.Range("D18").Value = "'8.47"
$ws.Range("E18").Value = "'  -1.97%  "
$ws.Range("D19").Value = "'3.325.86"
$ws.Range("E19").Value = "'  -2.47%  "
$ws.Range("D20").Value = "'0.526"
$ws.Range("E20").Value = "'  +2.82%  "
$ws.Range("D21").Value = "'16.92"
$ws.Range("E21").Value = "'  -2.88%  "
$ws.Range("D22").Value = "'501.50"
$ws.Range("E22").Value = "'  -0.57%  "
$ws.Range("B23").Value = "'SuiNetwork"
$ws.Range("C23").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").Value = "'3.36"
$ws.Range("E23").Value = "'  -2.11%  "
$ws.Range("B24").Value = "'Uniswap"
$ws.Range("C24").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").Value = "'10.42"
$ws.Range("E24").Value = "'  -4.19%  "
$ws.Range("E25").Value = "'  -4.11%  "
$ws.Range("D26").Value = "'6.50"
$ws.Range("E26").Value = "'  +6.11%  "
$ws.Range("D27").Value = "'95.55"
$ws.Range("E27").Value = "'  -3.26%  "
$ws.Range("D28").Value = "'11.96"
$ws.Range("E28").Value = "'  -5.79%  "
$ws.Range("B29").Value = "'WrappedeETH"
$ws.Range("C29").Value = "'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "'3.499.03"
$ws.Range("E29").Value = "'  -2.87%  "
$ws.Range("B30").Value = "'Hedera"
$ws.Range("C30").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.143"
$ws.Range("E30").Value = "'  -8.38%  "
$ws.Range("B31").Value = "'Dai"
$ws.Range("C31").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "'  +0.20%  "
$ws.Range("B32").Value = "'InternetComputer(DFINITY)"
$ws.Range("C32").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'10.88"
$ws.Range("E32").Value = "'  -4.28%  "
$ws.Range("B33").Value = "'Cronos"
$ws.Range("C33").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D33").Value = "'0.187"
$ws.Range("E33").Value = "'  -5.55%  "
$ws.Range("B34").Value = "'PancakeSwap"
$ws.Range("C34").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "'2.45"
$ws.Range("E34").Value = "'  +9.01%  "
$ws.Range("B35").Value = "'Binance-PegBSC-USD"
$ws.Range("C35").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "'  +0.27%  "
$ws.Range("B36").Value = "'PolygonEcosystemToken"
$ws.Range("C36").Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").Value = "'0.541"
$ws.Range("E36").Value = "'  -5.71%  "
$ws.Range("B37").Value = "'EthereumClassic"
$ws.Range("C37").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "'27.80"
$ws.Range("E37").Value = "'  -6.73%  "
$ws.Range("B38").Value = "'Fetch.AI"
$ws.Range("C38").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'1.45"
$ws.Range("E38").Value = "'  +3.42%  "
$ws.Range("B39").Value = "'RenderToken"
$ws.Range("C39").Value = "'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "'7.56"
$ws.Range("E39").Value = "'  -2.16%  "
$ws.Range("B40").Value = "'USDe"
$ws.Range("C40").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "'  +0.03%  "
$ws.Range("B41").Value = "'Kaspa"
$ws.Range("C41").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.150"
$ws.Range("E41").Value = "'  -2.33%  "
$ws.Range("B42").Value = "'Bittensor"
$ws.Range("C42").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'502.88"
$ws.Range("E42").Value = "'  -1.89%  "
$ws.Range("B43").Value = "'WhiteBITCoin"
$ws.Range("C43").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "'24.34"
$ws.Range("E43").Value = "'  -1.43%  "
$ws.Range("B44").Value = "'VeChain"
$ws.Range("C44").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0427"
$ws.Range("E44").Value = "'  +1.49%  "
$ws.Range("B45").Value = "'ARBITRUM"
$ws.Range("C45").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'0.824"
$ws.Range("E45").Value = "'  -3.61%  "
$ws.Range("B46").Value = "'MantraDAO"
$ws.Range("C46").Value = "'https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D46").Value = "'3.64"
$ws.Range("E46").Value = "'  -0.96%  "
$ws.Range("B47").Value = "'ImmutableX"
$ws.Range("C47").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").Value = "'1.65"
$ws.Range("E47").Value = "'  +5.48%  "
$ws.Range("B48").Value = "'Filecoin"
$ws.Range("C48").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "'5.44"
$ws.Range("E48").Value = "'  -0.19%  "
$ws.Range("B49").Value = "'Cosmos"
$ws.Range("C49").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'8.29"
$ws.Range("E49").Value = "'  +1.82%  "
$ws.Range("B50").Value = "'OKB"
$ws.Range("C50").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").Value = "'53.02"
$ws.Range("E50").Value = "'  +2.66%  "
$ws.Range("B51").Value = "'dogwifhat"
$ws.Range("C51").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'3.10"
$ws.Range("E51").Value = "'  -4.99%  "
